$wb = $excel.ActiveWorkbook

$newId = "d2562679-66db-4d8d-b08b-7655ded4db5e"
$newHash = "3a4de8036aab7508b025e0b5427f8e0059623ec4"

# --- Overview sheet: append new row for the new file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A3").Value = "$newId.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-17 12:27:15"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/$newId.md", "", "", "$newId.md")

# --- zh-cn sheet: append new row for the new file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A3").Value = "$newId.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-17 12:27:12"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/$newId.md", "", "", "$newId.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/$newId.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22e58273805f33c2677647823a2b5a483ccfe9e9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newId.$newHash.zh-cn.xlf", "", "", "$newId.$newHash.zh-cn.xlf")

# --- de-de sheet: append new row for the new file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A3").Value = "$newId.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-17 12:27:15"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/$newId.md", "", "", "$newId.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/$newId.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bfc791aa8658b10d0cdbbdadb83a60a04bac409/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newId.$newHash.de-de.xlf", "", "", "$newId.$newHash.de-de.xlf")

Write-Output "done"
